# Rename the embedded logo pictures' shape names:
#   - Pearson logo inline pictures (footers): image2.png -> image1.png
#   - BTEC logo inline picture (first-page header): image1.jpg -> image2.jpg
#
# InlineShape objects do not expose a settable .Name in the Word object
# model, so the rename is done the same way Word itself requires:
# convert the inline picture to a (floating) Shape, rename it there, then
# convert it back to an inline picture so the drawing stays wp:inline.

$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

function Rename-InlinePicture($inlineShape, [string]$newName) {
    $shape = $inlineShape.ConvertToShape()
    $shape.Name = $newName
    [void]$shape.ConvertToInlineShape()
}

# --- Footers: Pearson Edexcel logo (image2.png -> image1.png) ---
$footers = $sec.Footers
for ($i = 1; $i -le $footers.Count; $i++) {
    $footer = $footers.Item($i)
    if ($footer.Exists) {
        $shapes = $footer.Range.InlineShapes
        for ($j = 1; $j -le $shapes.Count; $j++) {
            $pic = $shapes.Item($j)
            if ($pic.AlternativeText -like "*PearsonLogo.png") {
                Rename-InlinePicture $pic "image1.png"
            }
        }
    }
}

# --- Headers: BTEC logo (image1.jpg -> image2.jpg) ---
$headers = $sec.Headers
for ($i = 1; $i -le $headers.Count; $i++) {
    $header = $headers.Item($i)
    if ($header.Exists) {
        $shapes = $header.Range.InlineShapes
        for ($j = 1; $j -le $shapes.Count; $j++) {
            $pic = $shapes.Item($j)
            if ($pic.AlternativeText -eq "BTec_Logo-Orange") {
                Rename-InlinePicture $pic "image2.jpg"
            }
        }
    }
}
